$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# O column: swap tied_teams order
$ws.Range("O52").Value = "['Colombia', 'Scotland']"

$o4 = "['Colombia', 'Scotland', 'Austria', 'Argentina']"
foreach ($r in 57,58,63,64,65,66,67,68,69,70,71,72,73) {
    $ws.Range("O$r").Value = $o4
}

$ws.Range("O101").Value = "['Saudi Arabia', 'United States']"
$ws.Range("O102").Value = "['Italy', 'Netherlands']"
$ws.Range("O104").Value = "['Netherlands', 'United States']"

# Rows 109-110: Bulgaria -> Argentina in Group D / top_four, plus change_flag/change_count updates
foreach ($r in 109,110) {
    $ws.Range("J$r").Value = "['Argentina', 6, 3, 6]"
    $ws.Range("M$r").Value = "['Argentina', 'Belgium', 'United States', 'Italy']"
    $ws.Range("Q$r").Value = 12
}

$ws.Range("P109").Value = 1
